$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 11: 2021年 ---
$ws.Range("A11").Value = "2021年"
$ws.Range("B11").Value = 315.7
$ws.Range("C11").Value = 57.3
$ws.Range("D11").Value = 11.5
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 310.3
$ws.Range("G11").Value = 623.6
$ws.Range("H11").Value = 79.8
$ws.Range("I11").Value = 292.5
$ws.Range("J11").Value = 37.7
$ws.Range("K11").Value = 49.6
$ws.Range("L11").Value = 36.2
$ws.Range("M11").Value = 67.8
$ws.Range("N11").Value = 50.8
$ws.Range("O11").Value = 222.1
$ws.Range("P11").Value = 76.5
$ws.Range("Q11").Value = 46.2
$ws.Range("R11").Value = 189.4
$ws.Range("S11").Value = 102.3
$ws.Range("T11").Value = 1020.4
$ws.Range("U11").Value = 7
$ws.Range("V11").Value = 560.9
$ws.Range("W11").Value = 90.90000000000001
$ws.Range("X11").Value = 3121.6
$ws.Range("Y11").Value = 625.9
$ws.Range("Z11").Value = 45.4
$ws.Range("AA11").Value = 347.7
$ws.Range("AB11").Value = 205.6
$ws.Range("AC11").Value = 119.8
$ws.Range("AD11").Value = 109.9
$ws.Range("AE11").Value = 11814.4
$ws.Range("AF11").Value = 1087
$ws.Range("AG11").Value = 314.7
$ws.Range("AH11").Value = 76
$ws.Range("AI11").Value = 81.2
$ws.Range("AJ11").Value = 22.8
$ws.Range("AK11").Value = 199.2
$ws.Range("AL11").Value = 306.3
$ws.Range("AM11").Value = 359.1
$ws.Range("AN11").Value = 42.3
$ws.Range("AO11").Value = 186.9
$ws.Range("AP11").Value = 198.6
$ws.Range("AQ11").Value = 116.1

# --- Row 12: 2022年 ---
$ws.Range("A12").Value = "2022年"
$ws.Range("B12").Value = 414.4
$ws.Range("C12").Value = 76.8
$ws.Range("D12").Value = 10.6
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 366.3
$ws.Range("G12").Value = 1139.7
$ws.Range("H12").Value = 169.4
$ws.Range("I12").Value = 418.6
$ws.Range("J12").Value = 50.5
$ws.Range("K12").Value = 61.3
$ws.Range("L12").Value = 69.2
$ws.Range("M12").Value = 45.9
$ws.Range("N12").Value = 56.1
$ws.Range("O12").Value = 388
$ws.Range("P12").Value = 74.7
$ws.Range("Q12").Value = 53.1
$ws.Range("R12").Value = 208.8
$ws.Range("S12").Value = 153.9
$ws.Range("T12").Value = 1223.9
$ws.Range("U12").Value = 2.2
$ws.Range("V12").Value = 480.2
$ws.Range("W12").Value = 206.3
$ws.Range("X12").Value = 2609.5
$ws.Range("Y12").Value = 645
$ws.Range("Z12").Value = 36.9
$ws.Range("AA12").Value = 1061.2
$ws.Range("AB12").Value = 142.8
$ws.Range("AC12").Value = 248.4
$ws.Range("AD12").Value = 96.59999999999999
$ws.Range("AE12").Value = 15568.1
$ws.Range("AF12").Value = 1756
$ws.Range("AG12").Value = 311.2
$ws.Range("AH12").Value = 152.5
$ws.Range("AI12").Value = 80.59999999999999
$ws.Range("AJ12").Value = 46.9
$ws.Range("AK12").Value = 295.8
$ws.Range("AL12").Value = 154.4
$ws.Range("AM12").Value = 680.5
$ws.Range("AN12").Value = 59.2
$ws.Range("AO12").Value = 152.3
$ws.Range("AP12").Value = 1300.9
$ws.Range("AQ12").Value = 67.5

# --- Copy style/format from row 10 (A10 label style) onto the new year labels ---
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A10").Copy()
$ws.Range("A12").PasteSpecial(-4122)
